$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# [name="Carol"]   "Only those who fear not sacrifice..." -> use single quotes instead of double quotes
$ws.Range("C25").Value = '[name="Carol"]   ''Only those who fear not sacrifice and possess the true, dauntless bloodline of Kazimierz may open the path.''
'

# [name="\"Captain\""]  You bunch... -> [name="'Captain'"]  You bunch...
$ws.Range("C51").Value = '[name="''Captain''"]  You bunch... How dare you show your faces here!
'

# [name="\"Captain\""]  Assemble! Everyone, get over here! -> [name="'Captain'"]  Assemble! Everyone, get over here!
$ws.Range("C52").Value = '[name="''Captain''"]  Assemble! Everyone, get over here!
'
